# Add a new trade record (row 4) to the BIIB random trade sheet,
# mirroring the formatting of the existing row 3 (date format on
# column A, boolean format on columns B/G/I, and the special style
# on columns A and G) and then filling in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 3's formatting down to row 4 so the new row keeps the same
# number formats/styles (e.g. style index on A/G) without introducing
# any new style definitions.
$ws.Range("A3:I3").Copy($ws.Range("A4:I4"))

# Overwrite with the new trade's values.
$ws.Range("A4").Value = 42633.676655092589
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = 10022
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 304.83
$ws.Range("F4").Value = 303.5
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = -0.44
$ws.Range("I4").Value = $false
